$d = $word.ActiveDocument

# The document contains a single table (AMOVA Phi results).
$t = $d.Tables.Item(1)

# Widen the first grid column (3456 -> 3737 twips == 172.8 -> 186.85 pt)
$t.Columns.Item(1).Width = 186.85

# Bump the height of the "body1" row (612 -> 617 twips == 30.6 -> 30.85 pt)
$t.Rows.Item(2).Height = 30.85

# Relabel "Populations" -> "Sampling sites" in the two cells that mention it.
$d.Content.Find.Execute("Populations:Total", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Sampling sites:Total", 2) | Out-Null

$d.Content.Find.Execute("Populations:Urban/rural groups", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Sampling sites:Urban/rural groups", 2) | Out-Null
